# Add two appendix tables to the workbook:
#   TableA1_discrepancy_ratios       - discrepancy ratio pivoted by country/year
#   TableA2_fire_loss_corr_sorted    - Table4_fire_loss_corr, sorted by Pearson correlation desc
#
# Headings use a clean "Country (ISO3)" / plain years style (no underscores),
# matching the header formatting already used on the other tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# TableA1_discrepancy_ratios
# ---------------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsA1 = $wb.Worksheets.Add($null, $lastSheet)
$wsA1.Name = "TableA1_discrepancy_ratios"

$wsA1.Cells.Item(1, 1).Value = "Country (ISO3)"
$wsA1.Cells.Item(1, 2).Value = 2016
$wsA1.Cells.Item(1, 3).Value = 2019
$wsA1.Cells.Item(1, 4).Value = 2023

$a1Data = @(
    @("ARG", 1.095996440146589, 0.8027770353969796, 1.201690072363169),
    @("BOL", 2.345918699311887, 3.65258765793092, 2.898758930285535),
    @("BRA", 2.260176146447278, 0.9505106100654231, 0.8420619971346361),
    @("CHL", 47.9543009184719, 9.413422224394401, 21.4823597782565),
    @("COL", 1.907426298753638, 1.390356715730083, 1.341511600125703),
    @("ECU", 0.6264072316900283, 0.464226115823015, 1.224966079563679),
    @("GUY", 2.582132641486334, 2.491279382962717, 3.556369363481102),
    @("PER", 1.128650839916678, 1.341478153156763, 0.7385605665069043),
    @("PRY", 1.111295844283485, 1.307075365266392, 1.102986604058451),
    @("SUR", 1.73008185564031, 1.039166524348505, 1.374601431569341),
    @("URY", $null, $null, $null),
    @("VEN", 1.51988724059989, 1.107808190268713, 1.17198458393493)
)

$r = 2
foreach ($row in $a1Data) {
    $wsA1.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne $null) { $wsA1.Cells.Item($r, 2).Value = $row[1] }
    if ($row[2] -ne $null) { $wsA1.Cells.Item($r, 3).Value = $row[2] }
    if ($row[3] -ne $null) { $wsA1.Cells.Item($r, 4).Value = $row[3] }
    $r++
}

# ---------------------------------------------------------------------------
# TableA2_fire_loss_corr_sorted
# ---------------------------------------------------------------------------

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsA2 = $wb.Worksheets.Add($null, $lastSheet2)
$wsA2.Name = "TableA2_fire_loss_corr_sorted"

$wsA2.Cells.Item(1, 1).Value = "Country (ISO3)"
$wsA2.Cells.Item(1, 2).Value = "Pearson correlation (Hansen vs VIIRS)"
$wsA2.Cells.Item(1, 3).Value = "Slope (ha per detection)"
$wsA2.Cells.Item(1, 4).Value = "Intercept (ha)"
$wsA2.Cells.Item(1, 5).Value = "p-value"

$a2Data = @(
    @("GUY", 0.8064161669000441, 5.418961741659722, 3310.845789785792, 0.008647339703521884),
    @("BOL", 0.7274446375198309, 2.231106557890573, 124092.3774577219, 0.02633763317206443),
    @("CHL", 0.6955870722410734, 421.964322851663, 75213.85905813953, 0.03746482915214577),
    @("PER", 0.6419039721847954, 2.341126097512947, 116768.3891580127, 0.0623488058854041),
    @("SUR", 0.6072982535445874, 4.996231459783445, 9119.974236294633, 0.08283402092999219),
    @("BRA", 0.5996083391776642, 4.483865718106943, -304395.6066071852, 0.08788451999342801),
    @("COL", 0.5703876575075997, 4.937424933639071, 41133.41418615848, 0.1087766912558362),
    @("ECU", 0.5523347760594717, 6.214116157655649, 15063.92295309163, 0.1230531024244793),
    @("VEN", 0.4367948232317739, 0.9928052043698148, 37126.2516346622, 0.239779237989031),
    @("ARG", 0.3513733988467822, 0.438411872807653, 172281.8475824944, 0.3538042294944216),
    @("URY", 0.1523869624050938, 67.02476141753287, 17366.58302568799, 0.6955047057618148),
    @("PRY", -0.3384923673147068, -25.27065980765401, 292937.1593391935, 0.3729219732774069)
)

$r = 2
foreach ($row in $a2Data) {
    $wsA2.Cells.Item($r, 1).Value = $row[0]
    $wsA2.Cells.Item($r, 2).Value = $row[1]
    $wsA2.Cells.Item($r, 3).Value = $row[2]
    $wsA2.Cells.Item($r, 4).Value = $row[3]
    $wsA2.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------------
# Header formatting - copy the bold/centered/bordered style already used for
# the other table headers (e.g. Table4's header row) onto the new headers.
# ---------------------------------------------------------------------------

$styleSrc = $wb.Worksheets.Item("Table4_fire_loss_corr").Range("A1")
$styleSrc.Copy()
$wsA1.Range("A1:D1").PasteSpecial(-4122)
$styleSrc.Copy()
$wsA2.Range("A1:E1").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Restore the original active sheet/selection (the workbook opened on the
# first sheet; adding sheets must not change that).
$wb.Worksheets.Item(1).Activate() | Out-Null
$wb.Worksheets.Item(1).Range("A1").Select() | Out-Null
